$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-9 (columns D, K, L, M, O, P change per-row;
# a brand-new row 9 is appended, row 8 is left untouched).
$rows = @(
    @{ Row = 2;  D = 44335; K = 18000; L = 20000; M = 19000; O = "Provincia de Limarí";  P = 760 },
    @{ Row = 3;  D = 44454; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí";  P = 540 },
    @{ Row = 4;  D = 44384; K = 12000; L = 13000; M = 12500; O = "Región de Coquimbo";   P = 500 },
    @{ Row = 5;  D = 44188; K = 18000; L = 20000; M = 19000; O = "Región Metropolitana"; P = 760 },
    @{ Row = 6;  D = 44316; K = 16000; L = 18000; M = 17000; O = "Región Metropolitana"; P = 680 },
    @{ Row = 7;  D = 44160; K = 9000;  L = 10000; M = 9500;  O = "Región Metropolitana"; P = 380 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($n, 11).Value = $r.K   # K: Precio mínimo
    $ws.Cells.Item($n, 12).Value = $r.L   # L: Precio máximo
    $ws.Cells.Item($n, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($n, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($n, 16).Value = $r.P   # P: Precio $/Kg
}

# New row 9: the record that used to live in row 4 before the shuffle.
$ws.Cells.Item(9, 1).Value  = 11
$ws.Cells.Item(9, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value  = "Bíobío"
$ws.Cells.Item(9, 4).Value  = 44162
$ws.Cells.Item(9, 5).Value  = 8
$ws.Cells.Item(9, 6).Value  = 100112026
$ws.Cells.Item(9, 7).Value  = "Haba"
$ws.Cells.Item(9, 8).Value  = "Sin especificar"
$ws.Cells.Item(9, 9).Value  = "Primera"
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 7500
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 7750
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Región Metropolitana"
$ws.Cells.Item(9, 16).Value = 310
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# Match the date-column style used by the other rows (D2:D8 -> style index 2 / numFmt 165).
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
